# Update the "Circumference" column (C) values with the recomputed
# circumference figures produced by the new circumference algorithm.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1412.285482287407
    3  = 3432.942671537399
    4  = 872.6559777259827
    5  = 1880.934319734573
    6  = 489.9726504087448
    7  = 1847.478476047516
    8  = 3115.232871055603
    9  = 843.7859205007553
    10 = 578.9259679317474
    11 = 1559.044929623604
    12 = 856.4549672603607
    13 = 1726.140382528305
    14 = 2997.989354014397
    15 = 690.7320977449417
    16 = 2322.757694482803
    17 = 1313.966721653938
    18 = 1710.667293906212
    19 = 565.6366448402405
    20 = 4006.544826507568
    21 = 2114.723204612732
    22 = 1070.989018917084
    23 = 1182.385987520218
    24 = 1089.366654634476
    25 = 597.5777697563171
    26 = 1151.089524149895
    27 = 1361.322060108185
    28 = 1044.923003554344
    29 = 704.9503531455994
    30 = 544.2396762371063
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
